$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.250.65"
$ws.Range('E2').Value = '  -2.73%  '
$ws.Range('D3').Value = "'1.933.82"
$ws.Range('E4').Value = '  +0.84%  '
$ws.Range('D5').Value = "'321.11"
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').Value = "'1.012"
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('D7').Value = "'0.4745"
$ws.Range('E7').Value = '  -4.57%  '
$ws.Range('D8').Value = "'0.4055"
$ws.Range('E8').Value = '  -3.58%  '
$ws.Range('D9').Value = "'53.57"
$ws.Range('D10').Value = "'0.08500"
$ws.Range('E10').Value = '  -7.90%  '
$ws.Range('D11').Value = "'1.051"
$ws.Range('E11').Value = '  -4.11%  '
$ws.Range('D12').Value = "'22.29"
$ws.Range('E12').Value = '  -2.32%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = "'7.519"
$ws.Range('E13').Value = '  -4.02%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = "'1.894.75"
$ws.Range('E14').Value = '  -4.43%  '
$ws.Range('D15').Value = "'6.121"
$ws.Range('E15').Value = '  -4.86%  '
$ws.Range('D17').Value = "'89.94"
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').Value = "'0.00001070"
$ws.Range('E18').Value = '  -2.62%  '
$ws.Range('D19').Value = "'0.06610"
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').Value = "'18.20"
$ws.Range('E20').Value = '  -5.38%  '
$ws.Range('D21').Value = "'1.012"
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').Value = "'5.796"
$ws.Range('E22').Value = '  -2.37%  '
$ws.Range('D23').Value = "'28.347.99"
$ws.Range('E23').Value = '  -2.48%  '
$ws.Range('D24').Value = "'11.43"
$ws.Range('E24').Value = '  -4.90%  '
$ws.Range('D25').Value = "'2.310"
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').Value = "'2.236.77"
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').Value = "'155.15"
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').Value = "'20.20"
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('D29').Value = "'2.169"
$ws.Range('E29').Value = '  -3.66%  '
$ws.Range('D30').Value = "'5.771"
$ws.Range('E30').Value = '  -8.01%  '
$ws.Range('D31').Value = "'123.96"
$ws.Range('E31').Value = '  -1.64%  '
$ws.Range('D32').Value = "'0.9828"
$ws.Range('E32').Value = '  -5.76%  '
$ws.Range('D33').Value = "'0.09608"
$ws.Range('E33').Value = '  -2.17%  '
$ws.Range('D34').Value = "'1.444"
$ws.Range('E34').Value = '  -5.28%  '
$ws.Range('D35').Value = "'3.668"
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').Value = "'5.590"
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('D37').Value = "'9.287"
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('E38').Value = '  -4.10%  '
$ws.Range('D39').Value = "'0.06183"
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').Value = "'1.241"
$ws.Range('E40').Value = '  -6.25%  '
$ws.Range('D41').Value = "'0.6201"
$ws.Range('E41').Value = '  -3.59%  '
$ws.Range('D42').Value = "'11.15"
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').Value = "'0.1909"
$ws.Range('E44').Value = '  -3.47%  '
$ws.Range('D45').Value = "'1.324"
$ws.Range('E45').Value = '  -0.36%  '
$ws.Range('D46').Value = "'0.5920"
$ws.Range('E46').Value = '  -4.80%  '
$ws.Range('D47').Value = "'12.80"
$ws.Range('E47').Value = '  -3.48%  '
$ws.Range('D48').Value = "'2.047"
$ws.Range('E48').Value = '  -6.91%  '
$ws.Range('D49').Value = "'3.396"
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('D50').Value = "'0.06785"
$ws.Range('E50').Value = '  -3.04%  '
$ws.Range('D51').Value = "'110.06"
$ws.Range('E51').Value = '  -1.87%  '
